# Minor reformatting of several slides (Records deck).
#
# Summary of changes:
#  - Slide 16: widen/shift the "lambda expression..." callout textbox,
#    reword its second line, nudge the diamond + connector, and bring
#    the diamond back in front of the connector (z-order).
#  - Slide 19: shift the "If the actual physical memory..." textbox right.
#  - Slide 21: add breathing room before the "type Rectangle" code block
#    and tighten the inline comment spacing on that line.
#  - Slide 25: reword the FieldExpr cast example to use "var" and update
#    the corresponding emit() call; nudge the bottom note box up.

$p = $ppt.ActivePresentation
$EMU_PER_PT = 12700

function EmuToPt($emu) {
    return $emu / $EMU_PER_PT
}

# ---------------------------------------------------------------------
# Slide 16
# ---------------------------------------------------------------------
$s16 = $p.Slides.Item(16)

# TextBox 5: reposition/resize and reword second line.
$textBox5 = $s16.Shapes.Item("TextBox 5")
$textBox5.Left  = EmuToPt 2971800
$textBox5.Width = EmuToPt 5069016

$tr16 = $textBox5.TextFrame.TextRange
$para2 = $tr16.Paragraphs(2, 1)
$sub = $tr16.Characters($para2.Start, $para2.Length)
$sub.Text = "lambda expression to compute record size."

# Connector: Elbow 8 - remove the horizontal flip and move/shrink it.
$connector = $s16.Shapes.Item("Connector: Elbow 8")
$connector.HorizontalFlip = 0
$connector.Left = EmuToPt 5506308
$connector.Width = EmuToPt 0

# Diamond 6 - nudge left and bring back in front of the connector.
$diamond = $s16.Shapes.Item("Diamond 6")
$diamond.Left = EmuToPt 5414868
$diamond.ZOrder(2)   # msoBringForward: swap with the connector behind it

# ---------------------------------------------------------------------
# Slide 19
# ---------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$textBox1 = $s19.Shapes.Item("TextBox 1")
$textBox1.Left = EmuToPt 1387475

# ---------------------------------------------------------------------
# Slide 21
# ---------------------------------------------------------------------
$s21 = $p.Slides.Item(21)
$content21 = $s21.Shapes.Item("Content Placeholder 2")
$tr21 = $content21.TextFrame.TextRange

# Blank paragraph right before "type Rectangle = record ..." gets extra
# space-before (0 -> 5 pts == spcPts val="500").
$blankPara = $tr21.Paragraphs(7, 1)
$blankPara.ParagraphFormat.SpaceBefore = 5

# Tighten the inline comment spacing on the "type Rectangle" line.
$rectPara = $tr21.Paragraphs(8, 1)
$rectSub = $tr21.Characters($rectPara.Start, $rectPara.Length)
$rectSub.Text = "type Rectangle = record        // fields are records"

# ---------------------------------------------------------------------
# Slide 25
# ---------------------------------------------------------------------
$s25 = $p.Slides.Item(25)
$content25 = $s25.Shapes.Item("Content Placeholder 2")
$tr25 = $content25.TextFrame.TextRange

# "FieldExpr fieldExpr = (FieldExpr) expr;" -> "var fieldExpr = (FieldExpr) expr;"
# Replace just the leading "FieldExpr " (first run + following space run)
# with "var " so the remaining runs (and their Consolas formatting) are
# left untouched.
$castPara = $tr25.Paragraphs(3, 1)
$castSub = $tr25.Characters($castPara.Start, 10)
$castSub.Text = "var "

# "expr.emit();" -> "fieldExpr.emit();"
$emitPara = $tr25.Paragraphs(5, 1)
$emitSub = $tr25.Characters($emitPara.Start, 9)
$emitSub.Text = "fieldExpr.emit"

# Bottom note textbox moves up slightly.
$textBox5_25 = $s25.Shapes.Item("TextBox 5")
$textBox5_25.Top = EmuToPt 5296251
